# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计", containing the
#    per-fund holding breakdown for the new quarter.
# 2. Insert a new row at the top of the "总计" (summary) sheet's data
#    table for the "2022-Q3" quarter, pushing the existing rows down.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" sheet, positioned right after 总计
# (i.e. right before the existing "2022-Q2" sheet).
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

# Header row (style copied from the summary sheet's header cell, which
# carries the bold/border "s=2" look used across every quarter sheet).
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# A-column index cells (style "s=2", same as the summary sheet's index
# column).
$summary.Range("A2").Copy()
$q3.Range("A2:A7").PasteSpecial(-4122)

$fundRows = @(
    @{ Row=2; Idx=0; Code="007139"; Name="富国民裕进取沪港深成长精选混合A"; Size="10.56"; Pos="88.09"; Ratio="5.36"; Value="0.5660"; Rank=7 },
    @{ Row=3; Idx=1; Code="011556"; Name="富国民裕进取沪港深成长精选混合C"; Size="2.46";  Pos="88.09"; Ratio="5.36"; Value="0.1319"; Rank=7 },
    @{ Row=4; Idx=2; Code="012586"; Name="南方港股创新视野一年持有混合A"; Size="2.16";  Pos="50.74"; Ratio="2.33"; Value="0.0503"; Rank=9 },
    @{ Row=5; Idx=3; Code="001942"; Name="前海开源沪港深汇鑫灵活配置混合A"; Size="0.17"; Pos="87.24"; Ratio="4.77"; Value="0.0081"; Rank=5 },
    @{ Row=6; Idx=4; Code="012587"; Name="南方港股创新视野一年持有混合C"; Size="0.19";  Pos="50.74"; Ratio="2.33"; Value="0.0044"; Rank=9 },
    @{ Row=7; Idx=5; Code="001943"; Name="前海开源沪港深汇鑫灵活配置混合C"; Size="0.09"; Pos="87.24"; Ratio="4.77"; Value="0.0043"; Rank=5 }
)

# Columns that must be stored as literal text (matches the source data,
# where fund codes / sizes / ratios are text, not numbers — e.g. to keep
# leading zeros in codes like "007139").
$textRange = $q3.Range("B2:G7")
$textRange.NumberFormat = "@"

foreach ($r in $fundRows) {
    $row = $r.Row
    $q3.Range("A$row").Value = $r.Idx
    $q3.Range("B$row").Value = $r.Code
    $q3.Range("C$row").Value = $r.Name
    $q3.Range("D$row").Value = $r.Size
    $q3.Range("E$row").Value = $r.Pos
    $q3.Range("F$row").Value = $r.Ratio
    $q3.Range("G$row").Value = $r.Value
    $q3.Range("H$row").Value = $r.Rank
}

# Drop the temporary "@" text format again so the cells end up with no
# explicit style, matching the rest of the workbook's plain data cells.
$textRange.ClearFormats()

# ---------------------------------------------------------------------
# Step 2: insert the "2022-Q3" row into the 总计 (summary) sheet.
# ---------------------------------------------------------------------
$summary.Rows("2:2").Insert()

# Inserting a row carries formatting down from the header row above, so
# strip that before writing the plain (unstyled) data cells.
$summary.Range("B2:D2").ClearFormats()

# A2 should keep the "s=2" index-column styling used by every other row
# in column A — restore it explicitly (Insert() drops it).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 6
$summary.Range("D2").Value = 0.76

# Keep the originally-active "总计" tab selected (matches source workbook,
# which had activeTab pointing at sheet 0 both before and after the edit).
$summary.Activate()
